$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.341583490371704
$ws.Range("B1").Value = 1.419857621192932
$ws.Range("C1").Value = 3.863321304321289
$ws.Range("D1").Value = 3.290687561035156
$ws.Range("E1").Value = 1.040981531143188
